# Regenerate save_data to use K (strike count) values recomputed from the
# raw odds data instead of the old "Strike#" based values. This updates the
# "K" column (column G) for every data row (rows 2-27) with the newly
# calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K column values, keyed by worksheet row number.
$newK = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 4
    12 = 0
    13 = 0
    14 = 3
    15 = 3
    16 = 1
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
